$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(16).Insert()
$ws.Range("B16:Q16").Clear()
$ws.Range("A16").Value = "R1Av0005"
$ws.Range("R16").Value = $null

$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value = "R1Bv0005"
$ws.Range("D22").Value = 40.08
$ws.Range("H22").Value = 320.64
$ws.Range("I22").Value = 40.08
$ws.Range("J22").Value = 320.64
$ws.Range("K22").Value = 320.64
$ws.Range("L22").Value = 320.64
$ws.Range("M22").Value = 320.64
$ws.Range("N22").Value = 320.64
$ws.Range("O22").Value = 320.64
$ws.Range("P22").Value = 40.08
$ws.Range("Q22").Value = 40.08
$ws.Range("R22").Value = 40.08

$ws.Range("E32").Select()
